{"js": "// Replace the opening paragraph's single run (\"List of the Webpages we\n// have found.\") with two runs \u2014 \"This is a list of all the websites we\n// have referenced in our slack \" and \"channel\" \u2014 with a grammar-check\n// proofErr (gramStart/gramEnd) bracketing the second run, matching the\n// author's edit.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstPara = paragraphs.items[0];\n\n// Grab the paragraph's own OOXML so we can read its real <w:p ...>\n// opening tag (paraId/rsid/etc.) and preserve it verbatim instead of\n// inventing new paragraph identity attributes.\nconst wholeRange = firstPara.getRange(Word.RangeLocation.whole);\nconst ooxmlResult = wholeRange.getOoxml();\nawait context.sync();\n\nconst existingXml = ooxmlResult.value || \"\";\nconst openTagMatch = /<w:p\\b[^>]*>/.exec(existingXml);\nconst pOpenTag = openTagMatch ? openTagMatch[0] : \"<w:p>\";\n\nconst newParagraphXml =\n  pOpenTag +\n  '<w:r><w:t xml:space=\"preserve\">This is a list of all the websites we have referenced in our slack </w:t></w:r>' +\n  '<w:proofErr w:type=\"gramStart\"/>' +\n  \"<w:r><w:t>channel</w:t></w:r>\" +\n  '<w:proofErr w:type=\"gramEnd\"/>' +\n  \"</w:p>\";\n\nconst flatOpcPackage =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n  \"<pkg:xmlData>\" +\n  '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n  \"</Relationships>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\" xmlns:mc=\"http://schemas.openxmlformats.org/markup-compatibility/2006\" mc:Ignorable=\"w14\">' +\n  \"<w:body>\" +\n  newParagraphXml +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\nwholeRange.insertOoxml(flatOpcPackage, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Replace the opening paragraph's single run (\"List of the Webpages we\n# have found.\") with two runs - \"This is a list of all the websites we\n# have referenced in our slack \" and \"channel\" - with a grammar-check\n# proofErr (gramStart/gramEnd) bracketing the second run, matching the\n# author's edit.\n\n$d = $word.ActiveDocument\n$p = $d.Paragraphs(1)\n$rng = $p.Range\n\n# Read back the paragraph's own OOXML so the real <w:p ...> opening tag\n# (paraId/rsid/etc.) can be preserved verbatim instead of inventing new\n# paragraph identity attributes.\n$existingXml = $rng.WordOpenXML\n$pOpenTag = '<w:p>'\nif ($existingXml -match '<w:p\\b[^>]*>') {\n    $pOpenTag = $matches[0]\n}\n\n$newParaXml = $pOpenTag +\n    '<w:r><w:t xml:space=\"preserve\">This is a list of all the websites we have referenced in our slack </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>channel</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '</w:p>'\n\n$flatOpc = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\" xmlns:mc=\"http://schemas.openxmlformats.org/markup-compatibility/2006\" mc:Ignorable=\"w14\">' +\n    '<w:body>' + $newParaXml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n\n$rng.InsertXML($flatOpc) | Out-Null\n"}
